# Replace the "observing campaign" date/constellation paragraph wherever
# it appears in the document (it is repeated once per language table).
# The new text is inserted as a single run with no explicit run
# formatting (the old multi-run, explicitly-fonted text is discarded).

$d = $word.ActiveDocument

$newText = "Informace v této příručce jsou určeny pro pozorovací kampaň probíhající od Perseus: 16.-25. Ledna, 7.-16. listopadu, 6.-15. prosince"

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.StartsWith("Informace v t")) {
        $start = $p.Range.Start
        $end = $p.Range.End

        # Clear all existing runs (and their formatting) in the paragraph,
        # then type the new sentence back in as plain, unformatted text
        # so the resulting run carries no <w:rPr/>.
        $r = $d.Range($start, $end)
        $r.Text = ""

        $r2 = $d.Range($start, $start)
        $r2.InsertAfter($newText)
    }
}
